# Advent of Code 2025, day 2 - add the next runtime data point to the
# RuntimesChart worksheet (the bar chart reads its series from
# RuntimesChart!$A$3:$A$20 / $B$3:$B$20, so appending a row is enough to
# extend the chart).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 0.0038892599999999999

# Match the author's final selection state (B4 is the new value cell).
$ws.Range("B4").Select() | Out-Null
